$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.83
$ws.Range("H2").Value = 3.25
$ws.Range("J2").Value = 2.6
$ws.Range("K2").Value = 1.95
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("W2").Value = 5.5
$ws.Range("Y2").Value = 9
$ws.Range("AA2").Value = 19
$ws.Range("AC2").Value = 6.5
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 23
$ws.Range("AP2").Value = 26
$ws.Range("AR2").Value = 67
$ws.Range("AS2").Value = 251
$ws.Range("AT2").Value = 2.38
$ws.Range("AU2").Value = 9.5
$ws.Range("BB2").Value = 401

# Row 3
$ws.Range("G3").Value = 3.8
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 2.1
$ws.Range("J3").Value = 4.5
$ws.Range("L3").Value = 2.88
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("W3").Value = 8.5
$ws.Range("Y3").Value = 15
$ws.Range("AA3").Value = 41
$ws.Range("AH3").Value = 8.5
$ws.Range("AJ3").Value = 19
$ws.Range("AO3").Value = 23
$ws.Range("AT3").Value = 2.38
$ws.Range("AX3").Value = 12
$ws.Range("AY3").Value = 26
$ws.Range("BD3").Value = 126

# Row 4
$ws.Range("G4").Value = 2.55
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 3.25
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.75
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 11
$ws.Range("Z4").Value = 26
$ws.Range("AG4").Value = 7.5
$ws.Range("AH4").Value = 13
$ws.Range("AW4").Value = 4.75
$ws.Range("AZ4").Value = 51
$ws.Range("BB4").Value = 251

# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("N5").Value = 7.5
$ws.Range("U5").Value = 1.95
$ws.Range("V5").Value = 1.8

# Row 6
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11

# Row 7
$ws.Range("G7").Value = 1.57
$ws.Range("I7").Value = 6
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("X7").Value = 7
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 19
$ws.Range("AG7").Value = 15
$ws.Range("AM7").Value = 900
$ws.Range("AW7").Value = 7
